# The "<id>p144r_1</id>" tag text for the p144r_1 transcription block was
# split across three separate runs in the paragraph:
#   run 1: "<id>"    (Courier New, color 7f6000, sz 18)
#   run 2: "p144r_1" (plain, color 000000)
#   run 3: "</id>"   (Courier New, color 7f6000, sz 18)
#
# The edit collapses those three runs into a single run containing the
# whole "<id>p144r_1</id>" string, taking on the formatting of the first
# run (Courier New / 7f6000 / sz 18). Doing a Find & Replace over that
# exact (document-unique) text span achieves this merge, since Word
# rewrites the matched range as one run using the formatting found at
# the start of the match.

$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("<id>p144r_1</id>", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "<id>p144r_1</id>", 2)
